# Rename the sheet from "UniformF-HW45.xpc" to "UniformF"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "UniformF"

# Append a new row (row 16) replicating the pattern of row 15,
# but with value 14 in column A and the "HexGrid-60degTilt5degRes" label in column B
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1

# Match formatting of row 15 (bold, centered, bordered cell in column A)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
